$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("G1").Value = "Discard Query Params"
$ws.Range("H1").Value = "Keep Query Params"

# New data cells (row 2 - /old-page rule)
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

# New data cells (row 3 - /legacy-section rule)
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $true
